$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40

$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"

$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 4).Value = 45121

$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = 100112013
$ws.Cells.Item($row, 7).Value = "Alcachofa"
$ws.Cells.Item($row, 8).Value = "Argentina(o)"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 120
$ws.Cells.Item($row, 11).Value = 16000
$ws.Cells.Item($row, 12).Value = 16000
$ws.Cells.Item($row, 13).Value = 16000
$ws.Cells.Item($row, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 16).Value = 320
$ws.Cells.Item($row, 17).Value = 50
$ws.Cells.Item($row, 18).Value = "Hortaliza"
